$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 8, shifting existing rows 8-16 down to 9-17.
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with the latest weekly price record.
$ws.Cells.Item(8, 1).Value = 8
$ws.Cells.Item(8, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(8, 3).Value = "Coquimbo"
$ws.Cells.Item(8, 4).Value = 45159
$ws.Cells.Item(8, 5).Value = 4
$ws.Cells.Item(8, 6).Value = 100114002
$ws.Cells.Item(8, 7).Value = "Camote"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 400
$ws.Cells.Item(8, 11).Value = 16000
$ws.Cells.Item(8, 12).Value = 17000
$ws.Cells.Item(8, 13).Value = 16500
$ws.Cells.Item(8, 14).Value = '$/malla 18 kilos'
$ws.Cells.Item(8, 15).Value = "Perú"
$ws.Cells.Item(8, 16).Value = 917
$ws.Cells.Item(8, 17).Value = 18
$ws.Cells.Item(8, 18).Value = "Hortaliza"
